$wb = $excel.ActiveWorkbook

# Remember the currently active sheet so we can restore selection/focus afterwards
$originalActiveSheet = $wb.ActiveSheet

# Add the new worksheet after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Combined Name Type"

# Header row (copy the existing header formatting from another sheet so the
# new header cells land on the same shared style as the rest of the workbook)
$newSheet.Range("A1").Value = "Events Name"
$newSheet.Range("B1").Value = "Event Type Name"

$wb.Worksheets.Item("Event Type Name List").Range("A1").Copy() | Out-Null
$newSheet.Range("A1:B1").PasteSpecial(-4122)

# Data rows: Event Name -> Event Type Name ("Other" for all of these)
$data = @(
    @("Hiatt Library: Spring 2025", "Other"),
    @("Possible Program: Spring 2025", "Other"),
    @("Rise Together, Register Together", "Other"),
    @("Moodle: Spring 2025", "Other"),
    @("LinkedIn Photo Booth Pop-up (Sherman)", "Other"),
    @("Career Closet: Spring 2025", "Other"),
    @("TypeFocus: Spring 2025", "Other"),
    @("Hiration: Spring 2025", "Other"),
    @("HWL Applying to Law School: Spring 2025", "Other"),
    @("INT 89: Spring 2025", "Other"),
    @("Rise Together: Spring 2025", "Other"),
    @("HWL Work, Career & Life: Spring 2025", "Other"),
    @("WOW: Spring 2025 ", "Other")
)

$row = 2
foreach ($pair in $data) {
    $newSheet.Cells.Item($row, 1).Value = $pair[0]
    $newSheet.Cells.Item($row, 2).Value = $pair[1]
    $row++
}

# Restore the originally active sheet/selection
$originalActiveSheet.Activate()
